# Applies the "Finished the core attack portion" changes to the Game Skills sheet:
#  - Astral Magics (row 13) gains a fight_time_out_mod_bonus_per_level (H13) of 0.001
#  - New "Criticality" skill row (row 22) is appended after Casting Accuracy (row 21)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game Skills")

# Astral Magics row - add fight_time_out_mod_bonus_per_level bonus
$ws.Range("H13").Value = 0.001

# New Criticality skill row
$ws.Range("A22").Value = "Criticality"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = "As you level this skill over time the skill bonus will grow to close to 100%. This bonus is used when you attack enemies to determine if you land a critical attack or not. This only applies to spells (healing and damage) and weapons."
$ws.Range("D22").Value = 999
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0.001
$ws.Range("N22").Value = 0
